$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = "44.462.42"
$ws.Range("E2").Value = "  +2.49%  "

# Row 3
$ws.Range("D3").Value = "2.368.11"
$ws.Range("E3").Value = "  -0.13%  "

# Row 4
$ws.Range("E4").Value = "  +0.08%  "

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "0.675"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +3.32%  "

# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "239.35"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +2.86%  "

# Row 7
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "73.23"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +5.92%  "

# Row 9
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.550"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +20.00%  "

# Row 10
$ws.Range("E10").Value = "  +6.96%  "

# Row 11
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "29.67"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +11.81%  "

# Row 12
$ws.Range("E12").Value = "  +1.97%  "

# Row 13
$ws.Range("D13").Value = "2.717.84"
$ws.Range("E13").Value = "  -0.25%  "

# Row 14
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "16.93"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +8.22%  "

# Row 15
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "6.75"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +7.55%  "

# Row 16
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.901"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +6.59%  "

# Row 17
$ws.Range("D17").Value = "2.366.39"
$ws.Range("E17").Value = "  -0.24%  "

# Row 18
$ws.Range("D18").Value = "44.356.58"
$ws.Range("E18").Value = "  +2.15%  "

# Row 19
$ws.Range("E19").Value = "  +4.07%  "

# Row 20
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "77.64"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +4.71%  "

# Row 21
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "6.48"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +3.26%  "

# Row 22
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "255.82"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +2.59%  "

# Row 23
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "0.999"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -0.11%  "

# Row 24
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "3.77"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -3.90%  "

# Row 25
$ws.Range("E25").Value = "  +2.74%  "

# Row 26
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "10.48"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +4.85%  "

# Row 27
$ws.Range("E27").Value = "  +1.04%  "

# Row 28
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "22.47"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +0.07%  "

# Row 29
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "174.38"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -0.57%  "

# Row 30
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "1.60"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +6.13%  "

# Row 32
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.133"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +5.27%  "

# Row 33
$ws.Range("B33").Value = "Filecoin"
$ws.Range("C33").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "5.23"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +4.44%  "

# Row 34
$ws.Range("B34").Value = "Hedera"
$ws.Range("C34").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.0742"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +7.15%  "

# Row 35
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "5.21"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +3.68%  "

# Row 36
$ws.Range("E36").Value = "  +7.88%  "

# Row 37
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "2.44"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -2.96%  "

# Row 38
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "6.53"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +0.28%  "

# Row 39
$ws.Range("E39").Value = "  +7.22%  "

# Row 40
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "19.99"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +10.47%  "

# Row 41
$ws.Range("E41").Value = "  +0.03%  "

# Row 42
$ws.Range("E42").Value = "  -1.44%  "

# Row 43
$ws.Range("E43").Value = "  +4.03%  "

# Row 44
$ws.Range("B44").Value = "ARBITRUM"
$ws.Range("C44").Value = "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "1.18"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +0.57%  "

# Row 45
$ws.Range("B45").Value = "Cronos"
$ws.Range("C45").Value = "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.0983"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +3.06%  "

# Row 46
$ws.Range("B46").Value = "Aave"
$ws.Range("C46").Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "98.79"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -0.50%  "

# Row 47
$ws.Range("E47").Value = "  +1.80%  "

# Row 48
$ws.Range("B48").Value = "Algorand"
$ws.Range("C48").Value = "https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.184"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +12.00%  "

# Row 49
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "2.36"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +4.03%  "

# Row 50
$ws.Range("D50").Value = "1.445.15"
$ws.Range("E50").Value = "  -0.16%  "

# Row 51
$ws.Range("D51").Value = "2.591.67"
$ws.Range("E51").Value = "  -0.13%  "
